# Refresh the cryptos price / 1h-volume columns (D, E) with the latest
# scraped values. D-column prices are forced to text (leading quote
# prefix) so things like "58.546.78" / "3.156.94" are not reinterpreted
# as numbers/dates by Excel. E-column values already carry the padding
# spaces from the source scrape, which also keeps them text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character used in the PEPE price (row 28)
$sub3 = [char]0x2083

$ws.Range("D2").Value = [string]::Concat("'", '58.546.78')
$ws.Range("E2").Value = '  -2.04%  '

$ws.Range("D3").Value = [string]::Concat("'", '3.156.94')
$ws.Range("E3").Value = '  -3.55%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = [string]::Concat("'", '527.69')
$ws.Range("E5").Value = '  -4.69%  '

$ws.Range("D6").Value = [string]::Concat("'", '135.44')
$ws.Range("E6").Value = '  -3.04%  '

# Row 7 (USDC) is unchanged in this update.

$ws.Range("D8").Value = [string]::Concat("'", '3.153.61')
$ws.Range("E8").Value = '  -3.67%  '

$ws.Range("D9").Value = [string]::Concat("'", '0.445')
$ws.Range("E9").Value = '  -3.88%  '

$ws.Range("D10").Value = [string]::Concat("'", '7.32')
$ws.Range("E10").Value = '  -6.39%  '

$ws.Range("E11").Value = '  -6.70%  '

$ws.Range("D12").Value = [string]::Concat("'", '0.380')
$ws.Range("E12").Value = '  -5.69%  '

$ws.Range("D13").Value = [string]::Concat("'", '3.690.33')
$ws.Range("E13").Value = '  -3.73%  '

$ws.Range("E14").Value = '  -0.99%  '

$ws.Range("D15").Value = [string]::Concat("'", '25.54')
$ws.Range("E15").Value = '  -3.55%  '

$ws.Range("D16").Value = [string]::Concat("'", '3.149.37')
$ws.Range("E16").Value = '  -3.68%  '

$ws.Range("D17").Value = [string]::Concat("'", '58.499.40')
$ws.Range("E17").Value = '  -2.32%  '

$ws.Range("D18").Value = [string]::Concat("'", '0.0000153')
$ws.Range("E18").Value = '  -5.70%  '

$ws.Range("D19").Value = [string]::Concat("'", '5.80')
$ws.Range("E19").Value = '  -4.06%  '

$ws.Range("D20").Value = [string]::Concat("'", '13.11')
$ws.Range("E20").Value = '  -3.48%  '

$ws.Range("D21").Value = [string]::Concat("'", '7.97')
$ws.Range("E21").Value = '  -5.66%  '

$ws.Range("D22").Value = [string]::Concat("'", '344.16')
$ws.Range("E22").Value = '  -7.11%  '

$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").Value = [string]::Concat("'", '0.512')
$ws.Range("E24").Value = '  -2.88%  '

$ws.Range("D25").Value = [string]::Concat("'", '67.25')
$ws.Range("E25").Value = '  -6.99%  '

$ws.Range("D26").Value = [string]::Concat("'", '3.272.41')
$ws.Range("E26").Value = '  -3.90%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("D28").Value = [string]::Concat('0.0', $sub3, '0958')
$ws.Range("E28").Value = '  -5.77%  '

$ws.Range("E29").Value = '  +0.92%  '

$ws.Range("E30").Value = '  -1.69%  '

$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("E32").Value = '  -6.57%  '

$ws.Range("E33").Value = '  -6.26%  '

$ws.Range("E34").Value = '  +2.57%  '

$ws.Range("D35").Value = [string]::Concat("'", '21.51')
$ws.Range("E35").Value = '  -3.99%  '

$ws.Range("D36").Value = [string]::Concat("'", '4.87')
$ws.Range("E36").Value = '  -2.62%  '

$ws.Range("D37").Value = [string]::Concat("'", '158.93')
$ws.Range("E37").Value = '  -4.48%  '

$ws.Range("D38").Value = [string]::Concat("'", '6.28')
$ws.Range("E38").Value = '  -4.30%  '

$ws.Range("E39").Value = '  -8.39%  '

$ws.Range("D40").Value = [string]::Concat("'", '0.0689')
$ws.Range("E40").Value = '  -4.30%  '

$ws.Range("D41").Value = [string]::Concat("'", '3.182.85')
$ws.Range("E41").Value = '  -3.63%  '

$ws.Range("E42").Value = '  -2.10%  '

$ws.Range("D43").Value = [string]::Concat("'", '24.08')
$ws.Range("E43").Value = '  -5.66%  '

$ws.Range("E44").Value = '  -0.61%  '

$ws.Range("D45").Value = [string]::Concat("'", '0.697')
$ws.Range("E45").Value = '  -6.20%  '

$ws.Range("D46").Value = [string]::Concat("'", '3.94')
$ws.Range("E46").Value = '  -3.50%  '

$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("E48").Value = '  -5.96%  '

$ws.Range("D49").Value = [string]::Concat("'", '2.285.77')
$ws.Range("E49").Value = '  -0.91%  '

$ws.Range("D50").Value = [string]::Concat("'", '6.17')
$ws.Range("E50").Value = '  -1.97%  '

$ws.Range("D51").Value = [string]::Concat("'", '20.80')
$ws.Range("E51").Value = '  -1.49%  '
